# Auto-generated edit script applying Seraph Profits market-data refresh
# (scheduled runner update) per the supplied OOXML diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 147.5
$ws.Cells.Item(9, 9).Value = 147.5
$ws.Cells.Item(9, 11).Value = 147.5
$ws.Cells.Item(9, 13).Value = 21.5
$ws.Cells.Item(29, 8).Value = 149.33333
$ws.Cells.Item(29, 9).Value = 149.33333
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = 447.99999
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = -166.99999
$ws.Cells.Item(29, 14).ClearContents()
$ws.Cells.Item(38, 8).Value = 458.625
$ws.Cells.Item(38, 10).Value = 0
$ws.Cells.Item(38, 12).Value = 0
$ws.Cells.Item(38, 14).ClearContents()
$ws.Cells.Item(48, 8).Value = 0
$ws.Cells.Item(48, 9).Value = 0
$ws.Cells.Item(48, 10).Value = 0
$ws.Cells.Item(48, 11).Value = 0
$ws.Cells.Item(48, 12).Value = 0
$ws.Cells.Item(48, 13).ClearContents()
$ws.Cells.Item(48, 14).ClearContents()
$ws.Cells.Item(56, 8).Value = 0
$ws.Cells.Item(56, 9).Value = 0
$ws.Cells.Item(56, 10).Value = 0
$ws.Cells.Item(56, 11).Value = 0
$ws.Cells.Item(56, 12).Value = 0
$ws.Cells.Item(56, 13).ClearContents()
$ws.Cells.Item(56, 14).ClearContents()
$ws.Cells.Item(58, 8).Value = 5023.75
$ws.Cells.Item(58, 10).Value = 6666.6665
$ws.Cells.Item(58, 12).Value = 19999.9995
$ws.Cells.Item(58, 14).Value = -20299.9995
$ws.Cells.Item(96, 8).Value = 71429530
$ws.Cells.Item(96, 9).Value = 71429530
$ws.Cells.Item(96, 10).Value = 0
$ws.Cells.Item(96, 11).Value = 214288590
$ws.Cells.Item(96, 12).Value = 0
$ws.Cells.Item(96, 13).Value = -214287217
$ws.Cells.Item(96, 14).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(56, 8).Value = 30000
$ws.Cells.Item(56, 10).Value = 30000
$ws.Cells.Item(56, 12).Value = 30000
$ws.Cells.Item(56, 14).Value = -31484
$ws.Cells.Item(74, 8).Value = 593.1070999999999
$ws.Cells.Item(74, 10).Value = 2500
$ws.Cells.Item(74, 12).Value = 2500
$ws.Cells.Item(74, 14).Value = -4248
$ws.Cells.Item(77, 8).Value = 593.1070999999999
$ws.Cells.Item(77, 10).Value = 2500
$ws.Cells.Item(77, 12).Value = 12500
$ws.Cells.Item(77, 14).Value = -21236
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).ClearContents()
$ws.Cells.Item(86, 14).ClearContents()
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).ClearContents()
$ws.Cells.Item(89, 14).ClearContents()
$ws.Cells.Item(97, 9).Value = 500.0909
$ws.Cells.Item(97, 10).Value = 897
$ws.Cells.Item(97, 11).Value = 500.0909
$ws.Cells.Item(97, 12).Value = 897
$ws.Cells.Item(97, 13).Value = -4.090899999999976
$ws.Cells.Item(97, 14).Value = -1889
$ws.Cells.Item(132, 8).Value = 4166.6665
$ws.Cells.Item(132, 9).Value = 2500
$ws.Cells.Item(132, 10).Value = 5000
$ws.Cells.Item(132, 11).Value = 7500
$ws.Cells.Item(132, 12).Value = 15000
$ws.Cells.Item(132, 13).Value = -4970
$ws.Cells.Item(132, 14).Value = -20060

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(33, 8).Value = 11374.25
$ws.Cells.Item(33, 9).Value = 2997
$ws.Cells.Item(33, 10).Value = 14166.667
$ws.Cells.Item(33, 11).Value = 2997
$ws.Cells.Item(33, 12).Value = 14166.667
$ws.Cells.Item(33, 13).Value = -2661
$ws.Cells.Item(33, 14).Value = -14838.667
$ws.Cells.Item(107, 8).Value = 1412.1428
$ws.Cells.Item(107, 9).Value = 1375.3889
$ws.Cells.Item(107, 11).Value = 1375.3889
$ws.Cells.Item(107, 13).Value = 544.6111000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 151.04
$ws.Cells.Item(7, 9).Value = 105.52381
$ws.Cells.Item(7, 10).Value = 390
$ws.Cells.Item(7, 11).Value = 105.52381
$ws.Cells.Item(7, 12).Value = 390
$ws.Cells.Item(7, 13).Value = 7.476190000000003
$ws.Cells.Item(7, 14).Value = -616
$ws.Cells.Item(31, 8).Value = 4553.1665
$ws.Cells.Item(31, 9).Value = 3164.75
$ws.Cells.Item(31, 11).Value = 3164.75
$ws.Cells.Item(31, 13).Value = -2869.75
$ws.Cells.Item(34, 8).Value = 4553.1665
$ws.Cells.Item(34, 9).Value = 3164.75
$ws.Cells.Item(34, 11).Value = 3164.75
$ws.Cells.Item(34, 13).Value = -2962.75
$ws.Cells.Item(56, 8).Value = 0
$ws.Cells.Item(56, 9).Value = 0
$ws.Cells.Item(56, 11).Value = 0
$ws.Cells.Item(56, 13).ClearContents()
$ws.Cells.Item(59, 8).Value = 48993.25
$ws.Cells.Item(59, 9).Value = 37997.5
$ws.Cells.Item(59, 10).Value = 59989
$ws.Cells.Item(59, 11).Value = 37997.5
$ws.Cells.Item(59, 12).Value = 59989
$ws.Cells.Item(59, 13).Value = -36852.5
$ws.Cells.Item(59, 14).Value = -62279
$ws.Cells.Item(74, 8).Value = 46711.855
$ws.Cells.Item(74, 10).Value = 46711.855
$ws.Cells.Item(74, 12).Value = 46711.855
$ws.Cells.Item(74, 14).Value = -48459.855
$ws.Cells.Item(77, 8).Value = 46711.855
$ws.Cells.Item(77, 10).Value = 46711.855
$ws.Cells.Item(77, 12).Value = 140135.565
$ws.Cells.Item(77, 14).Value = -148871.565
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 14).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 99500
$ws.Cells.Item(37, 10).Value = 99500
$ws.Cells.Item(37, 12).Value = 298500
$ws.Cells.Item(37, 14).Value = -298724
$ws.Cells.Item(58, 8).Value = 674.5
$ws.Cells.Item(58, 9).Value = 450
$ws.Cells.Item(58, 11).Value = 1350
$ws.Cells.Item(58, 13).Value = -1222
$ws.Cells.Item(107, 8).Value = 53270.42
$ws.Cells.Item(107, 9).Value = 620
$ws.Cells.Item(107, 11).Value = 1860
$ws.Cells.Item(107, 13).Value = 60
$ws.Cells.Item(138, 8).Value = 4559.6665
$ws.Cells.Item(138, 10).Value = 4784.5
$ws.Cells.Item(138, 12).Value = 14353.5
$ws.Cells.Item(138, 14).Value = -24633.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(32, 8).Value = 15000
$ws.Cells.Item(32, 10).Value = 15000
$ws.Cells.Item(32, 12).Value = 15000
$ws.Cells.Item(32, 14).Value = -15592
$ws.Cells.Item(123, 8).Value = 34969.6
$ws.Cells.Item(123, 10).Value = 34987
$ws.Cells.Item(123, 12).Value = 34987
$ws.Cells.Item(123, 14).Value = -39887
$ws.Cells.Item(132, 8).Value = 4000
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 4000
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 12000
$ws.Cells.Item(132, 13).ClearContents()
$ws.Cells.Item(132, 14).Value = -17060

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1001.5
$ws.Cells.Item(7, 9).Value = 1001.5
$ws.Cells.Item(7, 11).Value = 1001.5
$ws.Cells.Item(7, 13).Value = -889.5
$ws.Cells.Item(41, 8).Value = 24499.75
$ws.Cells.Item(41, 9).Value = 24000
$ws.Cells.Item(41, 10).Value = 24666.334
$ws.Cells.Item(41, 11).Value = 24000
$ws.Cells.Item(41, 12).Value = 24666.334
$ws.Cells.Item(41, 13).Value = -23562
$ws.Cells.Item(41, 14).Value = -25542.334
$ws.Cells.Item(126, 8).Value = 1001.5
$ws.Cells.Item(126, 9).Value = 1001.5
$ws.Cells.Item(126, 11).Value = 3004.5
$ws.Cells.Item(126, 13).Value = -534.5
$ws.Cells.Item(132, 8).Value = 1698.4
$ws.Cells.Item(132, 10).Value = 4500
$ws.Cells.Item(132, 12).Value = 13500
$ws.Cells.Item(132, 14).Value = -18560
$ws.Cells.Item(134, 8).Value = 63000
$ws.Cells.Item(134, 10).Value = 71250
$ws.Cells.Item(134, 12).Value = 71250
$ws.Cells.Item(134, 14).Value = -81390

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 4999.5
$ws.Cells.Item(126, 9).Value = 4999.5
$ws.Cells.Item(126, 11).Value = 14998.5
$ws.Cells.Item(126, 13).Value = -12528.5
